$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (20) of data, continuing the series present in rows 2-19.
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -0.7200474048664085
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -2.055952042396259

# Column A uses a date/time style in the rest of the series - copy it down
# from the row above so the new cell matches (style index "2" in the xlsx).
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
